{"js": "// Replace the ten three-digit-by-one-digit multiplication prompts that\n// changed in this revision. Every \"old\" prompt string in the document is\n// unique, so a plain (non-wildcard) search-and-replace on each one is\n// sufficient and keeps the original run formatting (font/size) intact.\nconst replacements = [\n  [\"473\u00d73=\", \"383\u00d72=\"],\n  [\"781\u00d73=\", \"530\u00d78=\"],\n  [\"768\u00d77=\", \"705\u00d72=\"],\n  [\"999\u00d75=\", \"425\u00d78=\"],\n  [\"170\u00d76=\", \"146\u00d76=\"],\n  [\"441\u00d78=\", \"296\u00d73=\"],\n  [\"822\u00d79=\", \"333\u00d79=\"],\n  [\"363\u00d77=\", \"382\u00d75=\"],\n  [\"517\u00d79=\", \"607\u00d77=\"],\n  [\"561\u00d73=\", \"661\u00d76=\"],\n  [\"313\u00d78=\", \"883\u00d72=\"],\n  [\"134\u00d72=\", \"424\u00d78=\"],\n  [\"781\u00d77=\", \"341\u00d72=\"],\n  [\"680\u00d78=\", \"326\u00d78=\"],\n  [\"829\u00d74=\", \"259\u00d79=\"],\n  [\"225\u00d75=\", \"448\u00d79=\"],\n  [\"857\u00d79=\", \"308\u00d77=\"],\n  [\"798\u00d79=\", \"503\u00d79=\"],\n  [\"271\u00d75=\", \"133\u00d77=\"],\n  [\"671\u00d75=\", \"782\u00d78=\"],\n  [\"268\u00d78=\", \"925\u00d74=\"],\n  [\"653\u00d74=\", \"220\u00d74=\"],\n  [\"987\u00d77=\", \"308\u00d79=\"],\n  [\"929\u00d78=\", \"292\u00d79=\"],\n  [\"323\u00d78=\", \"118\u00d74=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the ten three-digit-by-one-digit multiplication prompts that\n# changed in this revision. Every \"old\" prompt string in the document is\n# unique, so Find/Replace (wdReplaceAll) on each one is sufficient and\n# leaves the surrounding run formatting (font/size) untouched.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"473\u00d73=\", \"383\u00d72=\"),\n    @(\"781\u00d73=\", \"530\u00d78=\"),\n    @(\"768\u00d77=\", \"705\u00d72=\"),\n    @(\"999\u00d75=\", \"425\u00d78=\"),\n    @(\"170\u00d76=\", \"146\u00d76=\"),\n    @(\"441\u00d78=\", \"296\u00d73=\"),\n    @(\"822\u00d79=\", \"333\u00d79=\"),\n    @(\"363\u00d77=\", \"382\u00d75=\"),\n    @(\"517\u00d79=\", \"607\u00d77=\"),\n    @(\"561\u00d73=\", \"661\u00d76=\"),\n    @(\"313\u00d78=\", \"883\u00d72=\"),\n    @(\"134\u00d72=\", \"424\u00d78=\"),\n    @(\"781\u00d77=\", \"341\u00d72=\"),\n    @(\"680\u00d78=\", \"326\u00d78=\"),\n    @(\"829\u00d74=\", \"259\u00d79=\"),\n    @(\"225\u00d75=\", \"448\u00d79=\"),\n    @(\"857\u00d79=\", \"308\u00d77=\"),\n    @(\"798\u00d79=\", \"503\u00d79=\"),\n    @(\"271\u00d75=\", \"133\u00d77=\"),\n    @(\"671\u00d75=\", \"782\u00d78=\"),\n    @(\"268\u00d78=\", \"925\u00d74=\"),\n    @(\"653\u00d74=\", \"220\u00d74=\"),\n    @(\"987\u00d77=\", \"308\u00d79=\"),\n    @(\"929\u00d78=\", \"292\u00d79=\"),\n    @(\"323\u00d78=\", \"118\u00d74=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
